$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "Teléfono" column header (was "Telétono").
$ws.Range("E1").Value = "Teléfono"

# Update the saved view state: zoom level and the last selected cell.
$excel.ActiveWindow.Zoom = 136
$ws.Range("D18").Select()
